# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Updates the numeric values of column G (header "K") for rows 2-30 on Sheet1
# to reflect the recalculated per-game strikeout totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 2
    6  = 2
    7  = 6
    8  = 8
    9  = 0
    10 = 3
    11 = 4
    12 = 5
    13 = 9
    14 = 1
    15 = 4
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 3
    24 = 6
    25 = 4
    26 = 6
    27 = 3
    28 = 3
    29 = 0
    30 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
